$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE2").Value = 33.3
$ws.Range("AG2").Value = 67.8
$ws.Range("AE6").Value = 24.4
$ws.Range("AE8").Value = 28.9
$ws.Range("AG8").Value = 17.8
$ws.Range("AE9").Value = 41.1
$ws.Range("AG9").Value = 11.1
$ws.Range("AG10").Value = 16.7
$ws.Range("AE12").Value = 70
$ws.Range("AG12").Value = 70
$ws.Range("AE13").Value = 36.7
$ws.Range("AG13").Value = 87.8
$ws.Range("AE14").Value = 6.7
$ws.Range("AG14").Value = 5
$ws.Range("AE15").Value = 91.09999999999999
$ws.Range("AE16").Value = 32.2
$ws.Range("AG17").Value = 88.90000000000001
$ws.Range("AE18").Value = 20
$ws.Range("AE19").Value = 5
$ws.Range("AE20").Value = 44.4
$ws.Range("AE21").Value = 11.1
$ws.Range("AG23").Value = 54.4
$ws.Range("AE24").Value = 75.59999999999999
$ws.Range("AE25").Value = 84.40000000000001
$ws.Range("AG26").Value = 74.40000000000001
$ws.Range("AE27").Value = 23.3
$ws.Range("AE28").Value = 95
$ws.Range("AG28").Value = 72.2
$ws.Range("AE29").Value = 15.6
$ws.Range("AE30").Value = 48.9
$ws.Range("AE31").Value = 53.3
$ws.Range("AE32").Value = 80
$ws.Range("AG32").Value = 91.09999999999999
$ws.Range("AE33").Value = 11.1
$ws.Range("AG33").Value = 5
$ws.Range("AE34").Value = 71.09999999999999
$ws.Range("AG34").Value = 64.40000000000001
$ws.Range("AE35").Value = 17.8
$ws.Range("AG35").Value = 10
$ws.Range("AE36").Value = 20
$ws.Range("AE37").Value = 46.7
$ws.Range("AG38").Value = 78.90000000000001
$ws.Range("AE39").Value = 62.2
$ws.Range("AG39").Value = 81.09999999999999
$ws.Range("AG40").Value = 95
$ws.Range("AE41").Value = 75.59999999999999
$ws.Range("AE42").Value = 56.7
$ws.Range("AG42").Value = 56.7
$ws.Range("AE43").Value = 40
$ws.Range("AE45").Value = 57.8
$ws.Range("AE46").Value = 87.8
$ws.Range("AG46").Value = 85.59999999999999
